$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3rs NF")
$r = $ws.Range("H24")
$r.Interior.Color = 255
Write-Host ("after red: pattern=" + $r.Interior.Pattern)
$r.Interior.Color = -4142
Write-Host ("after none: pattern=" + $r.Interior.Pattern + " color=" + $r.Interior.Color)
